$d = $word.ActiveDocument

# --- Locate the anchor paragraph: "... After I got the backend working and
# finished the video, I went to sleep." is the last paragraph of the
# 13.01.2022 entry. We append three new paragraphs right after it (and right
# before the pre-existing trailing empty paragraph / sectPr).
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $par = $d.Paragraphs.Item($i)
    if ($par.Range.Text -like "*After I got the backend working and finished the video, I went to sleep.*") {
        $anchor = $par
    }
}

$rng = $anchor.Range
$rng.MoveEnd(1, -1)          # exclude the paragraph's own end-of-paragraph mark
$rng.Collapse(0)              # collapse to just after the last visible character

# Insert a placeholder blank-paragraph marker, the date paragraph, and the
# first sentence of the new entry, all in one go (keeps the anchor
# paragraph's text untouched). The placeholder ("ZZZ") is needed because a
# paragraph that is created empty from the start keeps a phantom empty run;
# deleting real placeholder text afterwards produces a genuinely run-less
# empty <w:p>, matching how the rest of the document's blank paragraphs look.
$rng.InsertAfter("`rZZZ`r24.01.2022`rI continued the MEAN-stack module. The video ")

$blankPar = $d.Paragraphs.Item(94)
$br = $blankPar.Range
$br.MoveEnd(1, -1)
$br.Text = ""

# --- Append the remaining three sentences as separate runs (matching the
# source, which keeps each sentence in its own <w:r>). A plain InsertAfter
# would merge same-format text into the previous run, so after inserting we
# briefly flip Bold on/off on just the new text -- this forces the engine to
# keep it as a distinct run. We also re-assert LanguageID explicitly since
# runs created this way can otherwise lose the inherited "en-US" language.
function Add-Sentence([string]$text) {
    $last = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
    $r = $last.Range
    $r.MoveEnd(1, -1)
    $r.Collapse(0)
    $r.InsertAfter($text)

    $find = $d.Content
    $find.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $find.LanguageID = "en-US"
    $find.Bold = 1
    $find.Bold = 0
}

Add-Sentence "covered setting up Angular to the full-stack application. I have only done this with React, but not with Angular."
Add-Sentence " I decided to stick with the newest version of the Angular-cli since I found it very silly to use older, deprecated versions."
Add-Sentence " I had to resolve some issues related to bootstrap (the video used an older version). After resolving the issues I had no problem with the rest of the video."
